# "final changes to models" - update the AIC table's model-comparison
# block (rows 11:15) with refreshed fit statistics, and leave the
# workbook with the AIC sheet active/selected (cell N21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AIC")

# --- Row 11: Era + Size ---
$ws.Range("D11").Value = -226.60613511116901
$ws.Range("G11").Value = 0.40029706342872801
$ws.Range("H11").Value = 118.512857765375
$ws.Range("I11").Value = 0.40029706342872801

# --- Row 12: Size ---
$ws.Range("D12").Value = -226.29381678498001
$ws.Range("E12").Value = 0.31231832618911398
$ws.Range("F12").Value = 0.85542302775902102
$ws.Range("G12").Value = 0.34242332600124797
$ws.Range("H12").Value = 117.285797281379
$ws.Range("I12").Value = 0.74272038942997598

# --- Row 13: Era x Size ---
$ws.Range("D13").Value = -225.722048195772
$ws.Range("E13").Value = 0.88408691539712003
$ws.Range("F13").Value = 0.64272170364255499
$ws.Range("G13").Value = 0.25727961057002402
$ws.Range("H13").Value = 119.156798745774

# --- Row 14: Era ---
$ws.Range("D14").Value = -19.002887980029801
$ws.Range("E14").Value = 207.60324713113999
$ws.Range("F14").Value = "8.3085965868500602E-46"
$ws.Range("G14").Value = "3.32590681493003E-46"
$ws.Range("H14").Value = 13.6403328789038

# --- Row 15: Null model ---
$ws.Range("D15").Value = -13.3275398207001
$ws.Range("E15").Value = 213.27859529046901
$ws.Range("F15").Value = "4.8656568672486602E-47"
$ws.Range("G15").Value = "1.9477081556114701E-47"
$ws.Range("H15").Value = 9.7465285310396794

# Recalculate all dependent (ROUND) formulas against the refreshed inputs.
$excel.CalculateFull() | Out-Null

# Leave the AIC sheet as the active/selected tab with N21 selected,
# matching the window state captured at the time of the final save.
$ws.Activate() | Out-Null
$ws.Range("N21").Select() | Out-Null
